$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A33").Value = "Dumb Framebuffer (DFB)"
$ws.Range("A34").Value = "Cursor (RAMDAC_CU)"
$ws.Range("A34").Select()
